$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E6: was a hard-coded 0.5, now becomes a formula (6/19) with a custom date-like
# number format (numFmtId 165) applied, per the "questions_correct_by_level"
# weighting rework.
$ws.Range("E6").NumberFormat = "YYYY""년 ""M""월 ""D\일"
$ws.Range("E6").Formula = "=6/19"

# E7: weight reduced from 1 to 0.5 (plain value edit, no formula/style change)
$ws.Range("E7").Value = 0.5

# G11: new column added next to the F11 total, computing the remaining points
# to distribute across the other 10 categories.
$ws.Range("G11").Formula = "=(1000-F11)/10"

# Move the active selection, matching the author's final cursor position.
$ws.Range("F15").Select()
